$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Employee")

$ws.Range("A6").Value = "ROLE GROUP : RTGO Operator 2023-12-07T19:27:58.156908600"
$ws.Range("A7").Value = "92970163 - Glenna Lynch`nROLE : RTGO100 1701853905917"
$ws.Range("A8").Value = "ROLE GROUP : RTGO Operator 2023-12-08T06:46:06.587"
